$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 314
$ws.Range("J17").Value = 314
$ws.Range("L17").Value = 942
$ws.Range("N17").Value = -1278
$ws.Range("H19").Value = 1604.7084
$ws.Range("I19").Value = 977.7
$ws.Range("J19").Value = 2052.5715
$ws.Range("K19").Value = 977.7
$ws.Range("L19").Value = 2052.5715
$ws.Range("M19").Value = -802.7
$ws.Range("N19").Value = -2402.5715
$ws.Range("H111").Value = 7465.3076
$ws.Range("I111").Value = 9941
$ws.Range("J111").Value = 4577
$ws.Range("K111").Value = 29823
$ws.Range("L111").Value = 13731
$ws.Range("M111").Value = -26756
$ws.Range("N111").Value = -19865
$ws.Range("H118").Value = 908.9231
$ws.Range("I118").Value = 847
$ws.Range("K118").Value = 2541
$ws.Range("M118").Value = -884
$ws.Range("H132").Value = 2411
$ws.Range("I132").Value = 2345.5264
$ws.Range("K132").Value = 7036.5792
$ws.Range("M132").Value = -4506.5792

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2084
$ws.Range("I2").Value = 1918.6923
$ws.Range("J2").Value = 2442.1667
$ws.Range("K2").Value = 1918.6923
$ws.Range("L2").Value = 2442.1667
$ws.Range("M2").Value = -1805.6923
$ws.Range("N2").Value = -2668.1667
$ws.Range("H61").Value = 1834.8684
$ws.Range("I61").Value = 1764.7715
$ws.Range("J61").Value = 2652.6667
$ws.Range("K61").Value = 1764.7715
$ws.Range("L61").Value = 2652.6667
$ws.Range("M61").Value = -1552.7715
$ws.Range("N61").Value = -3076.6667
$ws.Range("H74").Value = 55255.945
$ws.Range("I74").Value = 59870.5
$ws.Range("J74").Value = 2957.6667
$ws.Range("K74").Value = 59870.5
$ws.Range("L74").Value = 2957.6667
$ws.Range("M74").Value = -58996.5
$ws.Range("N74").Value = -4705.6667
$ws.Range("H77").Value = 55255.945
$ws.Range("I77").Value = 59870.5
$ws.Range("J77").Value = 2957.6667
$ws.Range("K77").Value = 299352.5
$ws.Range("L77").Value = 14788.3335
$ws.Range("M77").Value = -294984.5
$ws.Range("N77").Value = -23524.3335
$ws.Range("H88").Value = 1722.6086
$ws.Range("J88").Value = 1367.0769
$ws.Range("L88").Value = 1367.0769
$ws.Range("N88").Value = -2179.0769
$ws.Range("H91").Value = 1722.6086
$ws.Range("J91").Value = 1367.0769
$ws.Range("L91").Value = 1367.0769
$ws.Range("N91").Value = -4175.0769
$ws.Range("H102").Value = 935.9032
$ws.Range("I102").Value = 866.6896400000001
$ws.Range("K102").Value = 866.6896400000001
$ws.Range("M102").Value = 755.3103599999999
$ws.Range("H110").Value = 9854.6
$ws.Range("I110").Value = 9854.6
$ws.Range("K110").Value = 9854.6
$ws.Range("M110").Value = -7809.6
$ws.Range("H116").Value = 2084
$ws.Range("I116").Value = 1918.6923
$ws.Range("J116").Value = 2442.1667
$ws.Range("K116").Value = 1918.6923
$ws.Range("L116").Value = 2442.1667
$ws.Range("M116").Value = 375.3077000000001
$ws.Range("N116").Value = -7030.1667
$ws.Range("H132").Value = 2144.7917
$ws.Range("I132").Value = 2158
$ws.Range("K132").Value = 6474
$ws.Range("M132").Value = -3944
$ws.Range("H136").Value = 1834.8684
$ws.Range("I136").Value = 1764.7715
$ws.Range("J136").Value = 2652.6667
$ws.Range("K136").Value = 5294.3145
$ws.Range("L136").Value = 7958.000100000001
$ws.Range("M136").Value = -2744.3145
$ws.Range("N136").Value = -13058.0001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2084
$ws.Range("I3").Value = 1918.6923
$ws.Range("J3").Value = 2442.1667
$ws.Range("K3").Value = 1918.6923
$ws.Range("L3").Value = 2442.1667
$ws.Range("M3").Value = -1804.6923
$ws.Range("N3").Value = -2670.1667
$ws.Range("H94").Value = 1695.4166
$ws.Range("I94").Value = 1166
$ws.Range("K94").Value = 1166
$ws.Range("M94").Value = -715
$ws.Range("H107").Value = 1922.1818
$ws.Range("I107").Value = 1922.1818
$ws.Range("K107").Value = 1922.1818
$ws.Range("M107").Value = -2.181800000000067
$ws.Range("H134").Value = 4017.9814
$ws.Range("I134").Value = 4879.1714
$ws.Range("J134").Value = 2431.5789
$ws.Range("K134").Value = 14637.5142
$ws.Range("L134").Value = 7294.736699999999
$ws.Range("M134").Value = -12102.5142
$ws.Range("N134").Value = -12364.7367

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 2539.3333
$ws.Range("I6").Value = 3014.5715
$ws.Range("K6").Value = 3014.5715
$ws.Range("M6").Value = -2901.5715
$ws.Range("H7").Value = 310.22223
$ws.Range("I7").Value = 386.2857
$ws.Range("K7").Value = 386.2857
$ws.Range("M7").Value = -273.2857
$ws.Range("H17").Value = 3550
$ws.Range("I17").Value = 3550
$ws.Range("K17").Value = 3550
$ws.Range("M17").Value = -3376
$ws.Range("H25").Value = 5500
$ws.Range("I25").Value = 5000
$ws.Range("K25").Value = 5000
$ws.Range("M25").Value = -4826
$ws.Range("H41").Value = 0
$ws.Range("I41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("H107").Value = 4801.7715
$ws.Range("I107").Value = 1065.6154
$ws.Range("J107").Value = 7009.5
$ws.Range("K107").Value = 1065.6154
$ws.Range("L107").Value = 7009.5
$ws.Range("M107").Value = 854.3846000000001
$ws.Range("N107").Value = -10849.5
$ws.Range("H132").Value = 6592.0933
$ws.Range("I132").Value = 6157.0557
$ws.Range("K132").Value = 18471.1671
$ws.Range("M132").Value = -15941.1671

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 1429.85
$ws.Range("I39").Value = 799.93335
$ws.Range("J39").Value = 3319.6
$ws.Range("K39").Value = 2399.80005
$ws.Range("L39").Value = 9958.799999999999
$ws.Range("M39").Value = -2105.80005
$ws.Range("N39").Value = -10546.8
$ws.Range("H55").Value = 6294.5
$ws.Range("I55").Value = 3299.5
$ws.Range("J55").Value = 6793.6665
$ws.Range("K55").Value = 9898.5
$ws.Range("L55").Value = 20380.9995
$ws.Range("M55").Value = -9721.5
$ws.Range("N55").Value = -20734.9995

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 28675.871
$ws.Range("I132").Value = 38504
$ws.Range("K132").Value = 115512
$ws.Range("M132").Value = -112982
$ws.Range("H141").Value = 42250
$ws.Range("J141").Value = 42250
$ws.Range("L141").Value = 42250
$ws.Range("N141").Value = -52610

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2010.4445
$ws.Range("J46").Value = 3131.8333
$ws.Range("L46").Value = 3131.8333
$ws.Range("N46").Value = -3507.8333
$ws.Range("H82").Value = 2237.818
$ws.Range("I82").Value = 2145.6428
$ws.Range("J82").Value = 2399.125
$ws.Range("K82").Value = 2145.6428
$ws.Range("L82").Value = 2399.125
$ws.Range("M82").Value = -1784.6428
$ws.Range("N82").Value = -3121.125
$ws.Range("H85").Value = 2237.818
$ws.Range("I85").Value = 2145.6428
$ws.Range("J85").Value = 2399.125
$ws.Range("K85").Value = 2145.6428
$ws.Range("L85").Value = 2399.125
$ws.Range("M85").Value = -897.6428000000001
$ws.Range("N85").Value = -4895.125
$ws.Range("H132").Value = 4788.4443
$ws.Range("I132").Value = 4788.4443
$ws.Range("K132").Value = 14365.3329
$ws.Range("M132").Value = -11835.3329

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 435972.22
$ws.Range("I136").Value = 435972.22
$ws.Range("K136").Value = 1307916.66
$ws.Range("M136").Value = -1305366.66

$wb.Worksheets.Item("CRP").Range("M41").ClearContents()
